$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-6 from 2023-10-05 (45204) to 2023-10-08 (45207)
$ws.Range("C2:C6").Value = 45207
